# Append " (Changed main)" right after the existing sentence
# "This is a Microsoft word document." in the first paragraph, as three
# separate runs: " (", "Changed main", ")" -- matching how Word keeps
# freshly-typed text in its own run(s) rather than silently re-merging it
# into the pre-existing run.
$d = $word.ActiveDocument

# Find the end of the sentence we need to extend (robust to its exact
# location, instead of hard-coding a character offset).
$anchor = $d.Content
$anchor.Find.Execute("This is a Microsoft word document.", $true, $false, `
                      $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertionPoint = $anchor.End

# Turn on revision tracking just for the insertions so each InsertAfter
# call lands in its own run instead of being coalesced into the run it
# touches (which is what happens to plain, untracked insertions that
# share formatting with their neighbour).
$d.TrackRevisions = $true

$r1 = $d.Range($insertionPoint, $insertionPoint)
$r1.InsertAfter(" (")

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("Changed main")

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(")")

$d.TrackRevisions = $false

# Accept each recorded insertion individually (rather than
# Document.AcceptAllRevisions, which also forces an unrelated full
# re-layout pass) so we end up with plain runs and no tracked-change
# markup, while leaving the rest of the document byte-for-byte alone.
foreach ($revision in $d.Revisions) {
    $revision.Accept()
}
